# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Fri Jul  5 13:57:16 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "0.999" / "7.50" / "1.03" that must
# stay literal text (Excel would otherwise auto-convert them to numbers and
# silently drop meaningful trailing zeros). Force the whole data range to
# Text format before writing the new values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '55.587.79'
$ws.Range("E2").Value = '  -2.65%  '

$ws.Range("D3").Value = '2.963.02'
$ws.Range("E3").Value = '  -4.58%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '491.65'
$ws.Range("E5").Value = '  -5.16%  '

$ws.Range("D6").Value = '133.09'
$ws.Range("E6").Value = '  +0.77%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").Value = '2.957.45'
$ws.Range("E8").Value = '  -4.58%  '

$ws.Range("D9").Value = '0.419'
$ws.Range("E9").Value = '  -5.29%  '

$ws.Range("D10").Value = '7.17'
$ws.Range("E10").Value = '  +1.43%  '

$ws.Range("E11").Value = '  -5.76%  '

$ws.Range("D12").Value = '0.349'
$ws.Range("E12").Value = '  -7.69%  '

$ws.Range("E13").Value = '  +1.08%  '

$ws.Range("D14").Value = '3.469.28'
$ws.Range("E14").Value = '  -4.89%  '

$ws.Range("D15").Value = '24.77'
$ws.Range("E15").Value = '  -1.89%  '

$ws.Range("D16").Value = '55.460.17'
$ws.Range("E16").Value = '  -2.73%  '

$ws.Range("D17").Value = '2.961.99'
$ws.Range("E17").Value = '  -4.86%  '

$ws.Range("D18").Value = '0.0000140'
$ws.Range("E18").Value = '  -5.07%  '

$ws.Range("D19").Value = '5.67'
$ws.Range("E19").Value = '  -0.43%  '

$ws.Range("D20").Value = '12.17'
$ws.Range("E20").Value = '  -4.87%  '

$ws.Range("D21").Value = '7.50'
$ws.Range("E21").Value = '  -4.63%  '

$ws.Range("D22").Value = '318.33'
$ws.Range("E22").Value = '  -7.10%  '

$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.28%  '

$ws.Range("D24").Value = '0.463'
$ws.Range("E24").Value = '  -7.16%  '

$ws.Range("D25").Value = '60.30'
$ws.Range("E25").Value = '  -11.19%  '

$ws.Range("D26").Value = '1.03'
$ws.Range("E26").Value = '  +2.58%  '

$ws.Range("E27").Value = '  +1.29%  '

$ws.Range("E28").Value = '  +0.05%  '

$ws.Range("D29").Value = '0.0₃0845'
$ws.Range("E29").Value = '  -8.26%  '

$ws.Range("D30").Value = '6.55'
$ws.Range("E30").Value = '  -1.18%  '

$ws.Range("D31").Value = '6.54'
$ws.Range("E31").Value = '  -4.37%  '

$ws.Range("D32").Value = '1.17'
$ws.Range("E32").Value = '  -1.18%  '

$ws.Range("D33").Value = '1.70'
$ws.Range("E33").Value = '  -7.56%  '

$ws.Range("D34").Value = '19.50'
$ws.Range("E34").Value = '  -8.90%  '

$ws.Range("D35").Value = '150.02'
$ws.Range("E35").Value = '  -3.77%  '

$ws.Range("D36").Value = '4.40'
$ws.Range("E36").Value = '  -7.42%  '

$ws.Range("E37").Value = '  -4.98%  '

$ws.Range("D38").Value = '5.69'
$ws.Range("E38").Value = '  -6.65%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.0651'
$ws.Range("E39").Value = '  -4.11%  '

$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = '23.22'
$ws.Range("E40").Value = '  -8.02%  '

$ws.Range("D41").Value = '2.994.34'
$ws.Range("E41").Value = '  -4.94%  '

$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  +0.51%  '

$ws.Range("D43").Value = '36.33'
$ws.Range("E43").Value = '  -9.53%  '

$ws.Range("E44").Value = '  -3.18%  '

$ws.Range("E45").Value = '  -6.70%  '

$ws.Range("D46").Value = '1.38'
$ws.Range("E46").Value = '  -4.74%  '

$ws.Range("E47").Value = '  -8.26%  '

$ws.Range("D48").Value = '2.134.76'
$ws.Range("E48").Value = '  -4.19%  '

$ws.Range("D49").Value = '0.0235'
$ws.Range("E49").Value = '  +1.47%  '

$ws.Range("D50").Value = '19.29'
$ws.Range("E50").Value = '  -2.63%  '

$ws.Range("D51").Value = '5.57'
$ws.Range("E51").Value = '  -8.47%  '
